# Change bonus targets from annual to quarterly percentages,
# and update the corresponding bonus target currency amounts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 (IC4) - Paige Duty
$ws.Range("K3").Value = 3.75
$ws.Range("M3").Value = 6750

# Row 4 (IC3) - Lee Latency
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = $null
$ws.Range("M4").Value = 4500

# Row 5 (IC3) - Mona Torr
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 90
$ws.Range("M5").Value = 4350

# Row 6 (IC2) - Robin Rollback
$ws.Range("K6").Value = 2.5
$ws.Range("L6").Value = 95
$ws.Range("M6").Value = 3000

# Row 7 (IC2) - Kenny Canary
$ws.Range("K7").Value = 2.5
$ws.Range("M7").Value = 2875

# Row 8 (IC3) - Tracey Loggins
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 115
$ws.Range("M8").Value = 4650

# Row 9 (IC3) - Sue Q. Ell
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 90
$ws.Range("M9").Value = 4440

# Row 10 (IC2) - Jason Blob
$ws.Range("K10").Value = 2.5
$ws.Range("M10").Value = 2950

# Row 11 (IC4) - Al Ert
$ws.Range("K11").Value = 3.75
$ws.Range("L11").Value = $null
$ws.Range("M11").Value = 6562.5

# Row 12 (IC3) - Addie Min
$ws.Range("K12").Value = 3
$ws.Range("L12").Value = 100
$ws.Range("M12").Value = 4560

# Row 13 (IC2) - Tim Out
$ws.Range("K13").Value = 2.5
$ws.Range("L13").Value = $null
$ws.Range("M13").Value = 2750

# Row 14 (IC3) - Barbie Que
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 4470
